$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage-string cell updates (Coin name, Link, Volume(1h) columns,
# and a few Price cells that already contain non-numeric-looking text such as
# thousand-grouped "43.645.45").
$plainUpdates = [ordered]@{
    "D2" = '43.645.45'
    "E2" = '  -0.15%  '
    "D3" = '2.321.88'
    "E3" = '  +5.10%  '
    "E4" = '  -0.01%  '
    "E5" = '  +0.65%  '
    "E6" = '  +6.46%  '
    "E7" = '  +1.35%  '
    "E8" = '  -0.02%  '
    "E9" = '  +2.20%  '
    "E10" = '  -3.45%  '
    "E11" = '  +1.46%  '
    "E12" = '  +6.21%  '
    "E13" = '  +0.16%  '
    "D14" = '2.665.89'
    "E14" = '  +4.83%  '
    "E15" = '  +4.39%  '
    "E16" = '  +8.61%  '
    "D17" = '2.331.61'
    "E17" = '  +4.72%  '
    "D18" = '43.649.06'
    "E19" = '  +1.79%  '
    "E20" = '  +4.66%  '
    "E21" = '  +1.40%  '
    "E22" = '  +2.92%  '
    "E23" = '  -4.13%  '
    "E24" = '  +8.21%  '
    "E25" = '  +0.05%  '
    "E26" = '  -6.54%  '
    "E27" = '  +3.26%  '
    "E28" = '  -4.94%  '
    "E29" = '  +0.90%  '
    "E30" = '  -0.84%  '
    "E31" = '  +11.23%  '
    "E32" = '  -1.37%  '
    "E33" = '  -0.38%  '
    "E34" = '  +1.82%  '
    "E35" = '  +1.31%  '
    "E36" = '  -0.19%  '
    "E37" = '  +0.90%  '
    "E38" = '  -2.80%  '
    "E39" = '  +1.73%  '
    "E40" = '  +11.39%  '
    "E41" = '  +15.36%  '
    "E42" = '  +20.02%  '
    "E43" = '  -1.16%  '
    "E44" = '  -1.19%  '
    "B45" = 'FraxShare'
    "C45" = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    "E45" = '  +7.78%  '
    "B46" = 'MultiversX'
    "C46" = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    "E46" = '  -6.11%  '
    "E47" = '  +2.99%  '
    "E48" = '  +0.35%  '
    "E49" = '  +0.38%  '
    "D50" = '2.546.91'
    "E50" = '  +4.85%  '
    "E51" = '  -3.23%  '
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Price cells that look like plain numbers ("267.56", "0.628", ...). These must stay
# TEXT (matching the original inlineStr cell type / preserving trailing zeros such as
# "70.90"), so prefix with a leading apostrophe to force a text entry, exactly as a
# user typing into a General-formatted cell would do to keep the literal digits.
$textForcedUpdates = [ordered]@{
    "D5" = '267.56'
    "D6" = '91.64'
    "D7" = '0.628'
    "D10" = '44.44'
    "D11" = '0.0931'
    "D12" = '7.98'
    "D15" = '15.23'
    "D16" = '0.851'
    "D20" = '6.25'
    "D21" = '70.90'
    "D22" = '238.39'
    "D26" = '2.47'
    "D27" = '11.13'
    "D28" = '3.37'
    "D30" = '38.83'
    "D32" = '172.43'
    "D33" = '0.0888'
    "D34" = '5.48'
    "D37" = '4.42'
    "D38" = '0.0346'
    "D40" = '2.32'
    "D41" = '0.233'
    "D43" = '12.10'
    "D44" = '5.39'
    "D45" = '8.97'
    "D46" = '60.66'
    "D48" = '100.42'
    "D49" = '1.19'
    "D51" = '0.426'
}
foreach ($addr in $textForcedUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $textForcedUpdates[$addr]
}
